$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New historical rows (1985-1995 backward extension) to be placed at the top,
# right after the header row, pushing the existing series down.
$newRows = @(
    @(31228,1985,1.058157908195101,1986,-0.2746413702905404),
    @(31593,1986,0.613400232493766,1987,-0.6878195673301057),
    @(31958,1987,-2.809251799599333,1988,-5.22296739181829),
    @(32324,1988,1.311815945521522,1989,-0.2360507986455929),
    @(32689,1989,4.646141329820397,1990,4.409069002718513),
    @(33054,1990,6.234545104864941,1991,7.271550582979214),
    @(33419,1991,9.12959361095953,1992,10.06345753894098),
    @(33785,1992,3.770303686471776,1993,4.833608027690683),
    @(34150,1993,-3.631379003451563,1994,-4.837330717035071),
    @(34515,1994,2.525367061038386,1995,2.890876292836841),
    @(34880,1995,1.86590761958525,1996,1.408414909230937)
)

# Original data rows (previously rows 2-31), to be shifted down by 11 rows.
$oldRows = @(
    @(35221,1996,-0.6785505155195604,1997,-1.188578696562748),
    @(35586,1997,1.695006817304967,1998,1.278023271204054),
    @(35950,1998,2.847469613938181,1999,2.968999233486636),
    @(36319,1999,0.8798915307165922,2000,0.9767994936995539),
    @(36676,2000,2.656958100102846,2001,2.731710176905877),
    @(37034,2001,1.426159487177681,2002,1.268063782307416),
    @(37399,2002,-0.09079305020236461,2003,0.1072756205915404),
    @(37756,2003,-0.2852992240488517,2004,-0.6434466511772374),
    @(38120,2004,1.25935754237172,2005,1.52755625946297),
    @(38484,2005,1.849854064762901,2006,2.608452643869552),
    @(38848,2006,1.134182589932542,2007,0.9849328351749564),
    @(39217,2007,3.046387543902274,2008,2.70919674139074),
    @(39583,2008,3.696216801135943,2009,4.482779680928051),
    @(39948,2009,-10.1884206506174,2010,-12.44886541640216),
    @(40310,2010,1.026818295716803,2011,0.6731564222281827),
    @(40676,2011,4.330112454756896,2012,4.541390901327569),
    @(41044,2012,1.087923448804795,2013,1.111295745068719),
    @(41409,2013,-0.6600930445666675,2014,-0.7631906159002333),
    @(41774,2014,2.439056157886133,2015,2.712995226103132),
    @(42137,2015,1.478929861176237,2016,1.676382600382564),
    @(42503,2016,1.880356773996161,2017,2.143034813277178),
    @(42867,2017,1.929223532179036,2018,2.187377386010447),
    @(43235,2018,1.919825538525988,2019,1.616739045460869),
    @(43600,2019,0.8478136117613833,2020,1.153059686387214),
    @(43966,2020,-4.237342230872454,2021,-5.891534789017571),
    @(44341,2021,-1.195442653096257,2022,-4.00599668311078),
    @(44706,2022,1.344786076589832,2023,0.1554161503046547),
    @(45071,2023,-1.149816025990236,2024,-1.603768089243041),
    @(45436,2024,-0.213677328227746,2025,-0.1073236571678571),
    @(45800,2025,0.4937964215598223,2026,0.8228551885782087)
)

# Write the new rows into rows 2..12 (these rows already carry the right
# formatting from the original sheet, so only values need to change).
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = 2 + $i
    $row = $newRows[$i]
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").Value = $row[4]
}

# Shift the original rows down by 11 (old row 2 -> row 13, ..., old row 31 -> row 42).
for ($i = 0; $i -lt $oldRows.Length; $i++) {
    $r = 13 + $i
    $row = $oldRows[$i]
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").Value = $row[4]
}

# Rows 32..42 are brand new (beyond the original A1:E31 dimension) so column A
# needs the same date-style formatting ("YYYY-MM-DD HH:MM:SS", bold, bordered,
# centered) already used throughout column A. Copy it from an existing,
# correctly-styled cell instead of rebuilding it property-by-property, so no
# stray/duplicate style records get created.
$ws.Range("A2").Copy()
$ws.Range("A32:A42").PasteSpecial(-4122)
$excel.CutCopyMode = 0
